$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 389, shifting existing rows 389-455 down to 390-456.
$ws.Rows("389").Insert()

# Populate the newly inserted row 389 with the new weekly price record.
$ws.Range("A389").Value = 4
$ws.Range("B389").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C389").Value = "Los Lagos"
$ws.Range("D389").Value = 45218
$ws.Range("E389").Value = 10
$ws.Range("F389").Value = "Fruta"
$ws.Range("G389").Value = 100108
$ws.Range("H389").Value = "Tropicales y subtropicales"
$ws.Range("I389").Value = 100108002
$ws.Range("J389").Value = "Mango"
$ws.Range("K389").Value = "Sin especificar"
$ws.Range("L389").Value = "Primera"
$ws.Range("M389").Value = 100
$ws.Range("N389").Value = 13000
$ws.Range("O389").Value = 13000
$ws.Range("P389").Value = 13000
$ws.Range("Q389").Value = "$/bandeja 4 kilos"
$ws.Range("R389").Value = "Brasil"
$ws.Range("S389").Value = 3250
$ws.Range("T389").Value = 4
